# Apply the "Add files via upload" edit:
#  - Strip the trailing ", FL" from every submarket name in column A
#    (two names also get shortened: "Sanibel Island" -> "Sanibel Is",
#    "Jacksonville Bay Meadows/Butler Blvd" -> "Jacksonville Bay Meadows")
#  - Widen column A to fit the (now differently-sized) labels
#  - Leave the sheet scrolled down with A1:A49 selected, matching the
#    view state captured in the saved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DailyMeanValues_Jan2017")

$submarkets = @{
    2  = "Boca Raton "
    3  = "Bradenton/Airport "
    4  = "Clearwater "
    5  = "Coral Springs/Pompano Beach "
    6  = "Daytona Area "
    7  = "Daytona Beach "
    8  = "Florida Central Area"
    9  = "Florida Panhandle Area"
    10 = "Fort Lauderdale/Beach "
    11 = "Fort Myers Beach/Sanibel Is "
    12 = "Fort Myers/Bonita Springs "
    13 = "Fort Pierce/Port St Lucie "
    14 = "Fort Walton Beach "
    15 = "Gainesville "
    16 = "Hollywood/Airport "
    17 = "Jacksonville Bay Meadows "
    18 = "Jacksonville Beaches "
    19 = "Jacksonville Other Areas "
    20 = "Jacksonville/Airport "
    21 = "Key West "
    22 = "Kissimmee East "
    23 = "Kissimmee West "
    24 = "Lake Buena Vista "
    25 = "Lakeland/Winter Haven "
    26 = "Melbourne/Palm Bay "
    27 = "Miami Airport/Civic Center "
    28 = "Miami Beach "
    29 = "Miami CBD/North "
    30 = "Miami South "
    31 = "Naples "
    32 = "Ocala "
    33 = "Orlando Central "
    34 = "Orlando International Drive "
    35 = "Orlando North "
    36 = "Orlando South "
    37 = "Panama City "
    38 = "Pensacola "
    39 = "Sarasota/Beaches "
    40 = "St Petersburg "
    41 = "Tallahassee "
    42 = "Tampa CBD/Airport "
    43 = "Tampa East "
    44 = "Tampa North/Busch Gardens "
    45 = "Tarpon Springs/North Shore "
    46 = "Titusville/Cocoa Beach "
    47 = "Upper Florida Keys"
    48 = "West Broward/Plantation "
    49 = "West Palm Beach "
}

foreach ($row in $submarkets.Keys) {
    $ws.Cells.Item($row, 1).Value = $submarkets[$row]
}

# Column A needs to be much wider now that it holds the un-suffixed names
$ws.Columns.Item(1).ColumnWidth = 53.140625

# Restore the saved view: scrolled so row 22 is at the top, whole column
# A (the data range) selected
$win = $wb.Windows.Item(1)
$ws.Range("A1:A49").Select()
$win.ScrollRow = 22
$win.ScrollColumn = 1
